# Trade #43 closed at 2026-02-17 15:30:14 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up figures to account for the
# newly closed trade, and appends the trade's row to both the "All Trades"
# log and its strategy-specific "MarketMaking" log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet roll-up figures
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.77   # Current Capital
$wsSummary.Range("B4").Value = 0.77      # Total P&L $
$wsSummary.Range("B5").Value = 0.36      # Total P&L %
$wsSummary.Range("B6").Value = 43        # Total Trades
$wsSummary.Range("B8").Value = 20        # Losing Trades
$wsSummary.Range("B9").Value = 32.56     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.77     # Capital
$wsStatus.Range("D4").Value = 43         # Trades
$wsStatus.Range("E4").Value = 0.77       # P&L $
$wsStatus.Range("F4").Value = 0.77       # P&L %
$wsStatus.Range("G4").Value = 32.56      # Win Rate %

# ---------------------------------------------------------------------
# Append the new trade (row 44) to the "All Trades" log.
# Column B holds a plain date-looking string ("2026-02-17"); assigning it
# straight would get auto-coerced into a date serial number by the COM
# layer's smart-entry logic, so it is entered with a leading apostrophe to
# force literal text (matching how the source data stores it), then the
# cell style is reset back to "Normal" so no stray "quote prefix"
# formatting is left behind.
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Cells.Item(44, 1).Value = 43
$wsAllTrades.Cells.Item(44, 2).Value = "'2026-02-17"
$wsAllTrades.Cells.Item(44, 2).Style = "Normal"
$wsAllTrades.Cells.Item(44, 3).Value = "15:30:07"
$wsAllTrades.Cells.Item(44, 4).Value = "MarketMaking"
$wsAllTrades.Cells.Item(44, 5).Value = "UP"
$wsAllTrades.Cells.Item(44, 6).Value = 0.62
$wsAllTrades.Cells.Item(44, 7).Value = 0.59
$wsAllTrades.Cells.Item(44, 8).Value = "CLOSED"
$wsAllTrades.Cells.Item(44, 9).Value = -4.8387
$wsAllTrades.Cells.Item(44, 10).Value = -0.03
$wsAllTrades.Cells.Item(44, 11).Value = 100.77
$wsAllTrades.Cells.Item(44, 12).Value = 0
$wsAllTrades.Cells.Item(44, 13).Value = 0
$wsAllTrades.Cells.Item(44, 14).Value = 0.6
$wsAllTrades.Cells.Item(44, 15).Value = "Normal spread capture: 19600 bps"
$wsAllTrades.Cells.Item(44, 16).Value = "early_exit"
$wsAllTrades.Cells.Item(44, 17).Value = 0.13

# ---------------------------------------------------------------------
# Same trade row, appended to the strategy-specific "MarketMaking" log.
# ---------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$wsMarketMaking.Cells.Item(44, 1).Value = 43
$wsMarketMaking.Cells.Item(44, 2).Value = "'2026-02-17"
$wsMarketMaking.Cells.Item(44, 2).Style = "Normal"
$wsMarketMaking.Cells.Item(44, 3).Value = "15:30:07"
$wsMarketMaking.Cells.Item(44, 4).Value = "MarketMaking"
$wsMarketMaking.Cells.Item(44, 5).Value = "UP"
$wsMarketMaking.Cells.Item(44, 6).Value = 0.62
$wsMarketMaking.Cells.Item(44, 7).Value = 0.59
$wsMarketMaking.Cells.Item(44, 8).Value = "CLOSED"
$wsMarketMaking.Cells.Item(44, 9).Value = -4.8387
$wsMarketMaking.Cells.Item(44, 10).Value = -0.03
$wsMarketMaking.Cells.Item(44, 11).Value = 100.77
$wsMarketMaking.Cells.Item(44, 12).Value = 0
$wsMarketMaking.Cells.Item(44, 13).Value = 0
$wsMarketMaking.Cells.Item(44, 14).Value = 0.6
$wsMarketMaking.Cells.Item(44, 15).Value = "Normal spread capture: 19600 bps"
$wsMarketMaking.Cells.Item(44, 16).Value = "early_exit"
$wsMarketMaking.Cells.Item(44, 17).Value = 0.13
